# Edits to w1761265_Report.docx per commit "updated the word doc".
# Uses Range.Text assignment (rather than Find.Execute's built-in Replace)
# so that straight apostrophes in the replacement text are not silently
# smart-quoted by autoformat-as-you-type.

$d = $word.ActiveDocument

function Replace-DocText($OldText, $NewText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($OldText)
    if (-not $found) {
        throw "Could not find text: $OldText"
    }
    $rng.Text = $NewText
}

$rsquo = [char]0x2019

# 1) "The data structure I choose is a LinkedList ... method." ->
#    "A LinkedList (queue) is created ... by making use of the poll() method
#     to remove and get the first element of the queue."
Replace-DocText `
    "The data structure I choose is a LinkedList (queue) that is created in the BFS (Breadth First Search) method." `
    "A LinkedList (queue) is created in the BFS (Breadth First Search) method by making use of the poll() method to remove and get the first element of the queue."

# 2) Queue / Reason for using Queues sentence simplification.
Replace-DocText `
    " Queue is a data structure with both ends open, indicating that one end is often used to enter data the other end is often used to exclude data. The reason for using Queues is due the searching or traversing algorithm used is BFS (Breadth First Search/Traversal)." `
    " Queue is a data structure with both ends open, one end is often used to enter data the other is used to exclude data. Reason for using Queues is due the searching or traversing algorithm used is BFS (Breadth First Search/Traversal)."

# 3) "BFS is an algorithm" -> "BFS is a searching algorithm"
Replace-DocText `
    "BFS is an algorithm which is used for traversing a graph and this uses queues to remember to capture the next vertex to start a search. " `
    "BFS is a searching algorithm which is used for traversing a graph and this uses queues to remember to capture the next vertex to start a search. "

# 4) "The reason why " -> "Reason why "
Replace-DocText `
    "The reason why " `
    "Reason why "

# 5) Augmenting-path sentence: add "(least number of edges)" and trailing
#    complexity clause.
Replace-DocText `
    (" for finding the augmenting path is that BFS promises to find the shortest possible path from source to sink where as DFS doesn" + $rsquo + "t. ") `
    (" for finding the augmenting path is that BFS promises to find the shortest possible (least number of edges) path from source to sink where as DFS doesn" + $rsquo + "t, this also reduces the worst-case time complexity. ")

# 6) Ford Fulkerson intro sentence rewrite.
Replace-DocText `
    "Ford Fulkerson is the algorithm that was used. In a given graph, the Ford-Fulkerson algorithm is used to find the maximal flow from the start vertex to the sink vertex. Any edge in a graph has a capacity." `
    "Algorithm used is Ford Fulkerson. Ford-Fulkerson algorithm is used to find the maximal flow from start vertex to sink vertex. Any edge in a graph has a capacity."

# 7) Source/Sink paragraph rewrite (adds the Greedy Algorithm sentence at
#    the end). Keep the straight apostrophe in "edge's" intact.
Replace-DocText `
    "Source and Sink are the two key vertices that are given to find the maximum flow between these vertices. The sink vertex will have all inward edges and no outward edges, while the root vertex will have all outward edges and no inward edges. There are also some important constraints to be followed which are the flow on an edge cannot exceed its maximum capacity of flow through that edge and except for the source and sink, any edge's incoming and outgoing flow would be equal." `
    "Source and Sink are the two key vertices to find the maximum flow between them. The sink vertex will have all inward edges and no outward edges, the source vertex will have all outward edges and no inward edges. The flow on an edge cannot exceed its maximum capacity of flow through that edge and except for the source and sink, any edge's incoming and outgoing flow would be equal. The Greedy Algorithm approach was considered in this case."

Write-Output "edits applied"
